$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Pure text replacements (safe to do first; they don't change paragraph
#    count, so paragraph indices used later remain unaffected).
# ---------------------------------------------------------------------------

$d.Content.Find.Execute(
    "here time is y=direction of releasing or participating according to",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "here time is z=direction of releasing or making participating according to",
    2) | Out-Null

$d.Content.Find.Execute(
    "e = m *n2 *lamda2/T2",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "e = (m  * n2  * lamda2) / T2",
    2) | Out-Null

$d.Content.Find.Execute(
    "c = n*lambda/T",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "c = n * lambda / T",
    2) | Out-Null

$d.Content.Find.Execute(
    "v = n/T",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "v = n / T",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Paragraph insertions. Done from the bottom of the document upward so
#    that paragraph index numbers determined by inspecting the *original*
#    document stay valid for each subsequent insertion point.
# ---------------------------------------------------------------------------

# The blank paragraph that used to separate "where n proportional to 1/T"
# from "frequency(v)..." (paragraph 17) now holds new text instead of being
# blank -- no new paragraph is inserted here, the existing blank one is
# filled in.
$pMicro = $d.Paragraphs.Item(17)
$pMicro.Range.Text = "on microscopic scale"

# After the blank paragraph (originally paragraph 7, right after the
# "...bonding" paragraph / before "e= mc2") -> new paragraph about the
# multiverse / Schrodinger's cat.
$pBlank = $d.Paragraphs.Item(7)
$pBlank.Range.InsertParagraphAfter()
$pMulti = $d.Paragraphs.Item(8)
$pMulti.Range.Text = "space is multiverse it can exists multiple ways example Schrodinger’s cat the probability shows  infinity or not present"
$pMulti.Range.InsertAfter(" if cat can exist or not.")

# After "then what is gravitational force which is bonding" (paragraph 6)
# -> two new paragraphs: the "earth floats..." paragraph and the
# "around it...dimension...space" paragraph (with a superscript "st").
$pBonding = $d.Paragraphs.Item(6)
$pBonding.Range.InsertParagraphAfter()
$pEarth = $d.Paragraphs.Item(7)
$pEarth.Range.Text = "earth floats due to vibrational energy taking its as big ball full mass which generates magnetic waves "

$pEarth.Range.InsertParagraphAfter()
$pDim = $d.Paragraphs.Item(8)
$pDim.Range.Text = "around it and where as in blackhole both electro magnetic waves which causes particles pull inside  and leave since, 1st dimension matter , second dimension as magnetic(gravitational waves) grid) , third dimension as time and forth dimension as  space"

# Make the "st" in "1st" superscript.
$dimText = $pDim.Range.Text
$stIndex = $dimText.IndexOf("1st") + 1
$stStart = $pDim.Range.Start + $stIndex
$stRange = $d.Range($stStart, $stStart + 2)
$stRange.Font.Superscript = $true

Write-Host "Final paragraph count: $($d.Paragraphs.Count)"
